# repull data, push all data, mean calculation
# Update column F (dSF) values for the affected rows to match the
# repulled/pushed data set.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4  = -1
    5  = -2
    6  = -1
    8  = 0
    9  = 2
    10 = -6
    11 = -1
    12 = -5
    13 = -1
    16 = -1
    18 = -2
    19 = 1
    20 = 2
    21 = 2
    22 = -2
    23 = 3
    24 = 1
    25 = 3
    28 = -1
    29 = -6
    30 = -5
    31 = 1
    33 = -1
    34 = -2
    35 = 4
    36 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
